$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing data rows 3-18 down to 4-19 to make room for a new row 2
$ws.Range("A3:E18").Copy()
$ws.Range("A4").PasteSpecial(-4122)

# Write out the refreshed forecast series (all rows recomputed)
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2007
$ws.Range("C2").Value = 11.13090654781819
$ws.Range("D2").Value = 2008
$ws.Range("E2").Value = 10.67037004222142

$ws.Range("A3").Value = 39765
$ws.Range("B3").Value = 2008
$ws.Range("C3").Value = 4.672550446571067
$ws.Range("D3").Value = 2009
$ws.Range("E3").Value = -0.7156496512470745

$ws.Range("A4").Value = 40130
$ws.Range("B4").Value = 2009
$ws.Range("C4").Value = -14.45332333832743
$ws.Range("D4").Value = 2010
$ws.Range("E4").Value = 7.857938327064184

$ws.Range("A5").Value = 40494
$ws.Range("B5").Value = 2010
$ws.Range("C5").Value = 8.600536527919633
$ws.Range("D5").Value = 2011
$ws.Range("E5").Value = 12.64892828543749

$ws.Range("A6").Value = 40862
$ws.Range("B6").Value = 2011
$ws.Range("C6").Value = 10.25770250047622
$ws.Range("D6").Value = 2012
$ws.Range("E6").Value = 10.40099841437159

$ws.Range("A7").Value = 41228
$ws.Range("B7").Value = 2012
$ws.Range("C7").Value = 4.639893381363169
$ws.Range("D7").Value = 2013
$ws.Range("E7").Value = 6.662398279632087

$ws.Range("A8").Value = 41592
$ws.Range("B8").Value = 2013
$ws.Range("C8").Value = 0.3058963467304165
$ws.Range("D8").Value = 2014
$ws.Range("E8").Value = 1.195213983078647

$ws.Range("A9").Value = 41957
$ws.Range("B9").Value = 2014
$ws.Range("C9").Value = 4.068173739091874
$ws.Range("D9").Value = 2015
$ws.Range("E9").Value = 7.055025120039615

$ws.Range("A10").Value = 42321
$ws.Range("B10").Value = 2015
$ws.Range("C10").Value = 4.984288257750213
$ws.Range("D10").Value = 2016
$ws.Range("E10").Value = 1.985659800779893

$ws.Range("A11").Value = 42689
$ws.Range("B11").Value = 2016
$ws.Range("C11").Value = 1.878184267712912
$ws.Range("D11").Value = 2017
$ws.Range("E11").Value = -0.3562142672005275

$ws.Range("A12").Value = 43053
$ws.Range("B12").Value = 2017
$ws.Range("C12").Value = 4.695933104194339
$ws.Range("D12").Value = 2018
$ws.Range("E12").Value = 6.493919935864634

$ws.Range("A13").Value = 43418
$ws.Range("B13").Value = 2018
$ws.Range("C13").Value = 4.892602738886098
$ws.Range("D13").Value = 2019
$ws.Range("E13").Value = -2.576675125869599

$ws.Range("A14").Value = 43783
$ws.Range("B14").Value = 2019
$ws.Range("C14").Value = 0.8049382522247184
$ws.Range("D14").Value = 2020
$ws.Range("E14").Value = 3.1919852842623

$ws.Range("A15").Value = 44159
$ws.Range("B15").Value = 2020
$ws.Range("C15").Value = -8.78417389973717
$ws.Range("D15").Value = 2021
$ws.Range("E15").Value = 6.942816049735523

$ws.Range("A16").Value = 44525
$ws.Range("B16").Value = 2021
$ws.Range("C16").Value = 5.110501195359984
$ws.Range("D16").Value = 2022
$ws.Range("E16").Value = 0.8094958705429534

$ws.Range("A17").Value = 44890
$ws.Range("B17").Value = 2022
$ws.Range("C17").Value = 5.120680133083599
$ws.Range("D17").Value = 2023
$ws.Range("E17").Value = 0.5542886326586061

$ws.Range("A18").Value = 45254
$ws.Range("B18").Value = 2023
$ws.Range("C18").Value = -0.5532735011319234
$ws.Range("D18").Value = 2024
$ws.Range("E18").Value = -3.561435976944571

$ws.Range("A19").Value = 45618
$ws.Range("B19").Value = 2024
$ws.Range("C19").Value = -1.069674659641462
$ws.Range("D19").Value = 2025
$ws.Range("E19").Value = 0.01743232028155184
